$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Sheet5")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "PieContig"
$ws.Range("A40:N151").EntireRow.Delete()
Write-Output $ws.UsedRange.Address()
